$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New test case row (row 24) describing the "ddl_023" scenario:
# verifying a dropped table is removed from information_schema.tables.
$ws.Cells.Item(24, 1).Value = "ddl_023"
$ws.Cells.Item(24, 2).Value = "y"
$ws.Cells.Item(24, 3).Value = "删除表后验证是否从information_schema.tables中移除"
$ws.Cells.Item(24, 4).Value = "DDL"
$ws.Cells.Item(24, 5).Value = "drop"
$ws.Cells.Item(24, 6).Value = "schema56"
$ws.Cells.Item(24, 8).Value = 'drop table $schema56'
$ws.Cells.Item(24, 9).Value = 'select `TABLE_CATALOG`,`TABLE_SCHEMA`,`TABLE_NAME`,`TABLE_TYPE`,`ENGINE`,`VERSION`,`ROW_FORMAT`,`TABLE_ROWS`,`AVG_ROW_LENGTH`,`DATA_LENGTH`,`MAX_DATA_LENGTH`,`INDEX_LENGTH`,`DATA_FREE`,`AUTO_INCREMENT`,`UPDATE_TIME`,`CHECK_TIME`,`TABLE_COLLATION`,`CHECKSUM`,`CREATE_OPTIONS`,`TABLE_COMMENT` from information_schema.tables where `TABLE_NAME`=''$schema56'' or `TABLE_SCHEMA` in (''MYSQL'', ''INFORMATION_SCHEMA'')'
$ws.Cells.Item(24, 10).Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_023.csv"
$ws.Cells.Item(24, 11).Value = "csv_containsAll"

# Match the author's final selection on the new row.
$ws.Range("J24").Select()
